$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SC 92" row (original row 28) and the "RM 232" row (original row 26).
# Delete the higher-numbered row first so the lower row index is unaffected by the shift.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Apply the individual cell value corrections (post row-deletion numbering).
$ws.Range("C3").Value = 11.2
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("E9").Value = -6.8
$ws.Range("E10").Value = -6.1
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("F14").Value = 17.76
$ws.Range("F19").Value = 17.81
$ws.Range("F20").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("C32").Value = 10.5
